$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.453.78'
$ws.Range('E2').Value = '  -1.10%  '
$ws.Range('D3').Value = '1.956.26'
$ws.Range('E3').Value = '  -3.96%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E5').Value = '  -2.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.616'
$ws.Range('E6').Value = '  -3.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.39'
$ws.Range('E7').Value = '  -7.92%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.372'
$ws.Range('E9').Value = '  -5.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '55.71'
$ws.Range('E10').Value = '  -4.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0833'
$ws.Range('E11').Value = '  +4.30%  '
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.92'
$ws.Range('E13').Value = '  -6.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.833'
$ws.Range('E14').Value = '  -8.29%  '
$ws.Range('D15').Value = '2.241.63'
$ws.Range('E15').Value = '  -4.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.58'
$ws.Range('E16').Value = '  -6.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.34'
$ws.Range('E17').Value = '  -3.68%  '
$ws.Range('D18').Value = '1.950.53'
$ws.Range('E18').Value = '  -4.36%  '
$ws.Range('D19').Value = '36.336.52'
$ws.Range('E19').Value = '  -1.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.14'
$ws.Range('E20').Value = '  -3.16%  '
$ws.Range('D21').Value = '0.0₃0877'
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.50'
$ws.Range('E22').Value = '  -2.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.06'
$ws.Range('E23').Value = '  -6.61%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.52'
$ws.Range('E25').Value = '  -0.63%  '
$ws.Range('E26').Value = '  -2.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.51'
$ws.Range('E27').Value = '  -3.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.46'
$ws.Range('E28').Value = '  +2.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.66'
$ws.Range('E29').Value = '  -3.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.122'
$ws.Range('E30').Value = '  -13.87%  '
$ws.Range('E31').Value = '  -2.29%  '
$ws.Range('E32').Value = '  -1.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.73'
$ws.Range('E33').Value = '  -7.10%  '
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.32'
$ws.Range('E35').Value = '  -4.70%  '
$ws.Range('B36').Value = 'THORChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.14'
$ws.Range('E36').Value = '  -4.28%  '
$ws.Range('B37').Value = 'BinanceUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('E39').Value = '  -9.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.91'
$ws.Range('E40').Value = '  -7.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0985'
$ws.Range('E41').Value = '  -0.72%  '
$ws.Range('E42').Value = '  -2.88%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.19'
$ws.Range('E43').Value = '  -5.41%  '
$ws.Range('E44').Value = '  -2.57%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.77'
$ws.Range('E45').Value = '  -8.66%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.04'
$ws.Range('E46').Value = '  -8.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.95'
$ws.Range('E47').Value = '  -5.57%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '1.346.40'
$ws.Range('E48').Value = '  -2.04%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.36'
$ws.Range('E49').Value = '  -4.94%  '
$ws.Range('E50').Value = '  -3.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '45.25'
$ws.Range('E51').Value = '  -0.79%  '
